$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value()
        $v2 = $c2.Value()
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# Columns B..AC correspond to columns 2..29
# Swap rows 122 and 123 (data columns only, keep column A/id fixed)
Swap-Rows 122 123 2 29

# Swap rows 148 and 149 (data columns only, keep column A/id fixed)
Swap-Rows 148 149 2 29
